$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1623
$ws.Range("I2").Value = 4436
$ws.Range("J2").Value = 17832
$ws.Range("K2").Value = 88
$ws.Range("L2").Value = 4920
$ws.Range("M2").Value = 298
$ws.Range("N2").Value = 3150
$ws.Range("P2").Value = 65
$ws.Range("Q2").Value = 31
$ws.Range("R2").Value = 247
$ws.Range("S2").Value = 1963
$ws.Range("T2").Value = 3172
$ws.Range("U2").Value = 242
$ws.Range("V2").Value = 28171
$ws.Range("W2").Value = 19
$ws.Range("X2").Value = 27901
$ws.Range("Y2").Value = 43
$ws.Range("Z2").Value = 440
